$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll = $wb.Worksheets.Item("全部类型")

# Column F ("想去人数" / want-to-go count) updates on the "展览" sheet
$exhibitionUpdates = @{
    3  = 403
    6  = 23
    7  = 255
    8  = 14184
    9  = 137
    10 = 105
    11 = 5695
    12 = 582
    13 = 62
    17 = 4
    19 = 172
    20 = 775
    23 = 10478
    24 = 1193
    26 = 70
    27 = 3722
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Same logical rows (offset by the extra row on this sheet) on the "全部类型" sheet
$allTypesUpdates = @{
    3  = 403
    7  = 23
    8  = 255
    9  = 14184
    10 = 137
    11 = 105
    12 = 5695
    13 = 582
    14 = 62
    18 = 4
    20 = 172
    21 = 775
    25 = 10478
    26 = 1193
    28 = 70
    29 = 3722
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}

$wb.Save()
